$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '21.767.22'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '1.541.60'
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = "'289.94"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = "'0.3918"
$ws.Range("E7").Value = '  +3.31%  '
$ws.Range("D8").Value = "'0.3184"
$ws.Range("E8").Value = '  -2.87%  '
$ws.Range("D9").Value = "'43.20"
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").Value = "'0.07195"
$ws.Range("E10").Value = '  -2.18%  '
$ws.Range("D11").Value = "'1.066"
$ws.Range("E11").Value = '  -6.15%  '
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = "'5.627"
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").Value = "'18.62"
$ws.Range("E14").Value = '  -6.37%  '
$ws.Range("D15").Value = "'6.605"
$ws.Range("E15").Value = '  -3.78%  '
$ws.Range("D16").Value = '1.548.13'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = "'0.00001108"
$ws.Range("E17").Value = '  +1.47%  '
$ws.Range("D18").Value = "'0.06574"
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").Value = "'83.12"
$ws.Range("E19").Value = '  -2.66%  '
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = "'6.151"
$ws.Range("E21").Value = '  -4.78%  '
$ws.Range("E22").Value = '  -4.55%  '
$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = '  -7.32%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'2.364"
$ws.Range("E24").Value = '  +4.04%  '
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").Value = '21.766.75'
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").Value = "'2.390"
$ws.Range("E26").Value = '  -5.30%  '
$ws.Range("D27").Value = "'145.10"
$ws.Range("E27").Value = '  -3.79%  '
$ws.Range("E28").Value = '  -3.28%  '
$ws.Range("D29").Value = "'4.838"
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").Value = '1.721.65'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = "'117.87"
$ws.Range("E31").Value = '  -2.71%  '
$ws.Range("D32").Value = "'0.9668"
$ws.Range("E32").Value = '  -13.65%  '
$ws.Range("D33").Value = "'5.900"
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("D34").Value = "'0.08207"
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").Value = "'8.970"
$ws.Range("E35").Value = '  -3.71%  '
$ws.Range("D36").Value = "'0.06128"
$ws.Range("E36").Value = '  -1.48%  '
$ws.Range("D37").Value = "'5.126"
$ws.Range("E37").Value = '  -2.79%  '
$ws.Range("D38").Value = "'0.02214"
$ws.Range("E38").Value = '  -3.93%  '
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("D40").Value = "'1.183"
$ws.Range("E40").Value = '  -4.05%  '
$ws.Range("E41").Value = '  -24.15%  '
$ws.Range("D42").Value = "'0.9998"
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").Value = "'10.64"
$ws.Range("E43").Value = '  -3.65%  '
$ws.Range("D44").Value = "'0.5770"
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("D45").Value = "'13.09"
$ws.Range("E45").Value = '  -4.88%  '
$ws.Range("D46").Value = "'3.738"
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").Value = "'0.5537"
$ws.Range("E47").Value = '  -4.17%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = "'117.55"
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'1.884"
$ws.Range("E49").Value = '  -5.10%  '
$ws.Range("E50").Value = '  -2.92%  '
$ws.Range("D51").Value = "'0.06734"
$ws.Range("E51").Value = '  -3.68%  '
